$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.671.27"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "3.283.99"
$ws.Range("E3").Value = "  +5.07%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.09"
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.53"
$ws.Range("E6").Value = "  +3.88%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.284.55"
$ws.Range("E8").Value = "  +5.14%  "
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("E10").Value = "  +3.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.42"
$ws.Range("E11").Value = "  +3.54%  "
$ws.Range("E12").Value = "  +3.23%  "
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("D15").Value = "3.824.14"
$ws.Range("E15").Value = "  +5.10%  "
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").Value = "3.280.92"
$ws.Range("E17").Value = "  +4.98%  "
$ws.Range("D18").Value = "63.724.55"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.83"
$ws.Range("E19").Value = "  +3.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.51"
$ws.Range("E20").Value = "  +1.98%  "
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("E22").Value = "  +4.57%  "
$ws.Range("E23").Value = "  +4.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.51"
$ws.Range("E24").Value = "  +4.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.18"
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.76"
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("E28").Value = "  +6.72%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  +2.93%  "
$ws.Range("E31").Value = "  +3.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.69"
$ws.Range("E32").Value = "  +8.04%  "
$ws.Range("E33").Value = "  -3.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.52"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  +3.36%  "
$ws.Range("E36").Value = "  +3.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.98"
$ws.Range("E37").Value = "  +1.94%  "
$ws.Range("D38").Value = "0.0₃0738"
$ws.Range("E38").Value = "  +9.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0396"
$ws.Range("E39").Value = "  +3.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "426.17"
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("D41").Value = "3.054.07"
$ws.Range("E41").Value = "  +4.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.32"
$ws.Range("E42").Value = "  +1.83%  "
$ws.Range("E43").Value = "  +1.99%  "
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.263"
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("E46").Value = "  +4.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.14"
$ws.Range("E47").Value = "  +3.40%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  +2.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.71"
$ws.Range("E50").Value = "  +3.54%  "
$ws.Range("E51").Value = "  +1.59%  "
